$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Card22")

# Update L30 date from 24/1/2025 to 24/1/2026
$ws.Range("L30").Value = "24/1/2026"

# Remove row 31 entirely (duplicate entry), shifting cells up
$ws.Rows.Item(31).Delete()
